$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newLabel = "Índice de volume de vendas no comércio varejista ampliado - Variação mensal (base: igual mês do ano anterior)"

$newValues = @{
    2  = 5.1
    3  = 2.9
    4  = -2.2
    5  = -11
    6  = -6.7
    7  = 6.9
    8  = 1.8
    9  = 4.1
    10 = 2.8
    11 = -2.7
    12 = 6.8
    13 = 3.611111111111111
    14 = -0.07777777777777777
    15 = -15.63333333333333
    16 = -4.211111111111111
    17 = 1.588888888888889
    18 = 1.122222222222222
    19 = 1.477777777777778
    20 = 6.366666666666667
    21 = -4.455555555555556
    22 = 3.7
    23 = 0.6
    24 = 0.1
    25 = -21.1
    26 = 2.6
    27 = -4
    28 = 4.7
    29 = -4.9
    30 = 8.300000000000001
    31 = -8.5
}

for ($row = 2; $row -le 31; $row++) {
    $ws.Cells.Item($row, 2).Value = $newLabel
    $ws.Cells.Item($row, 4).Value = $newValues[$row]
}
